$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0053929171847225616
$ws.Range("D2").Value = 0.062968923947104738
$ws.Range("E2").Value = 0.08410963464109325
$ws.Range("C3").Value = 0.0032024662792370711
$ws.Range("D3").Value = 0.037786355828070654
$ws.Range("E3").Value = 0.050340222523895299
$ws.Range("C4").Value = 0.0097579531301865124
$ws.Range("D4").Value = 0.18236023852174274
$ws.Range("E4").Value = 0.22061227485472515
$ws.Range("C5").Value = 0.0052340502429781442
$ws.Range("D5").Value = 0.085033182323522971
$ws.Range("E5").Value = 0.10555098432368948
$ws.Range("C6").Value = 0.010792808446874533
$ws.Range("D6").Value = 0.23412728283929227
$ws.Range("E6").Value = 0.27643604322548537
$ws.Range("C7").Value = 0.0073001574389876499
$ws.Range("D7").Value = 0.14790331044256422
$ws.Range("E7").Value = 0.17652038096154449
$ws.Range("C8").Value = 0.01171484875585397
$ws.Range("D8").Value = 0.32745813099486387
$ws.Range("E8").Value = 0.37338137066057386
$ws.Range("C9").Value = 0.0080009668851266282
$ws.Range("D9").Value = 0.18442183530854439
$ws.Range("E9").Value = 0.21578612237841971
$ws.Range("C10").Value = 0.011947485636022496
$ws.Range("D10").Value = 0.35897449010825483
$ws.Range("E10").Value = 0.40580968684876045
$ws.Range("C11").Value = 0.0089881976243184943
$ws.Range("D11").Value = 0.22260298408645685
$ws.Range("E11").Value = 0.25783727696347769
$ws.Range("C12").Value = 0.012684275237936939
$ws.Range("D12").Value = 0.36707566009007059
$ws.Range("E12").Value = 0.41679913701046428
$ws.Range("C13").Value = 0.010423061283285079
$ws.Range("D13").Value = 0.23849156376838099
$ws.Range("E13").Value = 0.27935061129719507
$ws.Range("C14").Value = 0.012225684717779425
$ws.Range("D14").Value = 0.28374372951408977
$ws.Range("E14").Value = 0.33166949117545441
$ws.Range("C15").Value = 0.010803144130047612
$ws.Range("D15").Value = 0.23706799892873737
$ws.Range("E15").Value = 0.27941699482096188
$ws.Range("C16").Value = 0.010862189644393845
$ws.Range("D16").Value = 0.19868698064363127
$ws.Range("E16").Value = 0.24126772143933506
$ws.Range("C17").Value = 0.013323199108590703
$ws.Range("D17").Value = 0.20945979331208917
$ws.Range("E17").Value = 0.26168756122195874
$ws.Range("C18").Value = 0.011598457333979141
$ws.Range("D18").Value = 0.061626797024620567
$ws.Range("E18").Value = 0.10709377205788106
$ws.Range("C19").Value = 0.014603580281820072
$ws.Range("D19").Value = 0.11745230299910191
$ws.Range("E19").Value = 0.17469924462292372
